$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 106, shifting existing rows 106:211 down to 107:212
$ws.Rows.Item(106).Insert()

# Populate the newly inserted row 106 with the new record
$ws.Range("A106").Value = 10
$ws.Range("B106").Value = "Vega Modelo de Temuco"
$ws.Range("C106").Value = "La Araucanía"
$ws.Range("D106").Value = 44781
$ws.Range("E106").Value = 9
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100104
$ws.Range("H106").Value = "Frutos de pepita"
$ws.Range("I106").Value = 100104003
$ws.Range("J106").Value = "Membrillo"
$ws.Range("K106").Value = "Champion"
$ws.Range("L106").Value = "Primera"
$ws.Range("M106").Value = 8
$ws.Range("N106").Value = 250000
$ws.Range("O106").Value = 250000
$ws.Range("P106").Value = 250000
$ws.Range("Q106").Value = "$/bins (450 kilos)"
$ws.Range("R106").Value = "Región de O'Higgins"
$ws.Range("S106").Value = 556
$ws.Range("T106").Value = 450
